$d = $word.ActiveDocument

# Update the date line at the top of the document
$d.Content.Find.Execute("2024-10-31 Thursday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-11-01 Friday", 2)

# Update the division-problem table. Addressing cells by (row, column)
# avoids any ambiguity from duplicate cell text (e.g. "10÷7=1, 3" occurs
# twice in the source document).
$tbl = $d.Tables.Item(1)

$tbl.Cell(1, 1).Range.Text  = "48÷4=12, 0"
$tbl.Cell(1, 2).Range.Text  = "41÷4=10, 1"
$tbl.Cell(1, 3).Range.Text  = "94÷7=13, 3"
$tbl.Cell(1, 4).Range.Text  = "54÷5=10, 4"
$tbl.Cell(1, 5).Range.Text  = "42÷4=10, 2"

$tbl.Cell(5, 1).Range.Text  = "64÷4=16, 0"
$tbl.Cell(5, 2).Range.Text  = "45÷9=5, 0"
$tbl.Cell(5, 3).Range.Text  = "77÷7=11, 0"
$tbl.Cell(5, 4).Range.Text  = "69÷7=9, 6"
$tbl.Cell(5, 5).Range.Text  = "79÷7=11, 2"

$tbl.Cell(9, 1).Range.Text  = "97÷7=13, 6"
$tbl.Cell(9, 2).Range.Text  = "72÷9=8, 0"
$tbl.Cell(9, 3).Range.Text  = "26÷2=13, 0"
$tbl.Cell(9, 4).Range.Text  = "94÷8=11, 6"
$tbl.Cell(9, 5).Range.Text  = "96÷8=12, 0"

$tbl.Cell(13, 1).Range.Text = "87÷4=21, 3"
$tbl.Cell(13, 2).Range.Text = "57÷9=6, 3"
$tbl.Cell(13, 3).Range.Text = "22÷9=2, 4"
$tbl.Cell(13, 4).Range.Text = "49÷6=8, 1"
$tbl.Cell(13, 5).Range.Text = "45÷7=6, 3"

$tbl.Cell(17, 1).Range.Text = "30÷7=4, 2"
$tbl.Cell(17, 2).Range.Text = "81÷9=9, 0"
$tbl.Cell(17, 3).Range.Text = "40÷2=20, 0"
$tbl.Cell(17, 4).Range.Text = "19÷5=3, 4"
$tbl.Cell(17, 5).Range.Text = "28÷7=4, 0"
